$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new log entry row (row 33)
$ws.Cells.Item(33, 1).Value = (Get-Date -Year 2025 -Month 1 -Day 9 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(33, 2).Value = "Upstairs mansion and implementing enemy"
$ws.Cells.Item(33, 3).Value = 4

# Match date formatting used by the rest of column A
$ws.Cells.Item(33, 1).NumberFormat = $ws.Cells.Item(32, 1).NumberFormat

# Update view state to match the authored workbook
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("B34").Select()
